$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "50.716.91"
$ws.Range("E2").Value = "  -1.17%  "

# Row 3
$ws.Range("D3").Value = "2.918.52"
$ws.Range("E3").Value = "  -1.55%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'374.64"
$ws.Range("E5").Value = "  -1.87%  "

# Row 6
$ws.Range("D6").Value = "'99.52"
$ws.Range("E6").Value = "  -3.07%  "

# Row 7
$ws.Range("E7").Value = "  -1.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").Value = "'0.574"
$ws.Range("E9").Value = "  -2.68%  "

# Row 10
$ws.Range("D10").Value = "'35.60"
$ws.Range("E10").Value = "  -2.52%  "

# Row 11
$ws.Range("E11").Value = "  -0.78%  "

# Row 12
$ws.Range("D12").Value = "'0.0844"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13
$ws.Range("D13").Value = "3.380.99"
$ws.Range("E13").Value = "  -1.80%  "

# Row 14
$ws.Range("D14").Value = "'17.95"
$ws.Range("E14").Value = "  -0.74%  "

# Row 15
$ws.Range("D15").Value = "'7.60"
$ws.Range("E15").Value = "  +1.77%  "

# Row 16
$ws.Range("E16").Value = "  +63.42%  "

# Row 17
$ws.Range("D17").Value = "2.922.83"
$ws.Range("E17").Value = "  -1.60%  "

# Row 18
$ws.Range("D18").Value = "'0.990"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19
$ws.Range("D19").Value = "50.706.50"
$ws.Range("E19").Value = "  -1.10%  "

# Row 20
$ws.Range("D20").Value = "'2.99"
$ws.Range("E20").Value = "  -7.16%  "

# Row 21
$ws.Range("D21").Value = "'12.22"
$ws.Range("E21").Value = "  -3.36%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -1.46%  "

# Row 23
$ws.Range("D23").Value = "'69.19"
$ws.Range("E23").Value = "  +0.72%  "

# Row 24
$ws.Range("D24").Value = "'265.36"
$ws.Range("E24").Value = "  +1.35%  "

# Row 25
$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  +9.14%  "

# Row 26
$ws.Range("E26").Value = "  -3.51%  "

# Row 27
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("D28").Value = "'7.08"
$ws.Range("E28").Value = "  -5.50%  "

# Row 29
$ws.Range("D29").Value = "'25.31"
$ws.Range("E29").Value = "  -1.69%  "

# Row 30
$ws.Range("E30").Value = "  -3.13%  "

# Row 31
$ws.Range("E31").Value = "  -4.58%  "

# Row 32
$ws.Range("E32").Value = "  +0.86%  "

# Row 33
$ws.Range("D33").Value = "'50.28"
$ws.Range("E33").Value = "  -1.24%  "

# Row 34
$ws.Range("E34").Value = "  -0.11%  "

# Row 35
$ws.Range("D35").Value = "'33.09"
$ws.Range("E35").Value = "  -2.76%  "

# Row 36
$ws.Range("D36").Value = "'0.0427"
$ws.Range("E36").Value = "  -3.94%  "

# Row 37
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("E38").Value = "  +2.26%  "

# Row 39
$ws.Range("E39").Value = "  -0.48%  "

# Row 40
$ws.Range("D40").Value = "'16.26"
$ws.Range("E40").Value = "  -4.38%  "

# Row 41
$ws.Range("E41").Value = "  -0.50%  "

# Row 42 (Stacks -> Monero)
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'122.83"
$ws.Range("E42").Value = "  +0.46%  "

# Row 43 (Monero -> Stacks)
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.40"
$ws.Range("E43").Value = "  -6.13%  "

# Row 44
$ws.Range("D44").Value = "'20.73"
$ws.Range("E44").Value = "  -2.03%  "

# Row 45
$ws.Range("E45").Value = "  -1.95%  "

# Row 46
$ws.Range("D46").Value = "'3.33"
$ws.Range("E46").Value = "  +2.87%  "

# Row 47
$ws.Range("E47").Value = "  -0.89%  "

# Row 48
$ws.Range("D48").Value = "1.990.40"
$ws.Range("E48").Value = "  -1.49%  "

# Row 49
$ws.Range("D49").Value = "'0.257"
$ws.Range("E49").Value = "  -6.28%  "

# Row 50
$ws.Range("D50").Value = "'0.0312"
$ws.Range("E50").Value = "  -6.83%  "

# Row 51
$ws.Range("E51").Value = "  +3.10%  "
